$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected with password "D382"; unprotect to edit values, then re-protect.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclaimer text (row 16, column A).
# Only the day portion of the date actually changed (2021-03-18 -> 2021-03-19),
# so locate and replace just that substring rather than rewriting the whole cell.
$disclaimerCell = $ws.Range("A16")
$disclaimerText = $disclaimerCell.Value()
$oldDate = "2021-03-18"
$newDate = "2021-03-19"
$dateIndex = $disclaimerText.IndexOf($oldDate)
$disclaimerCell.Characters($dateIndex + 1, $oldDate.Length).Text = $newDate

# Update Weight (D) and Percent Change (E) values for rows 2-13.
$ws.Range("D2").Value = 0.03054091421840216
$ws.Range("E2").Value = -0.02000000000000013
$ws.Range("D3").Value = 0.02363536082577141
$ws.Range("E3").Value = -0.004466501240694809
$ws.Range("D4").Value = 0.05149740948159683
$ws.Range("E4").Value = -0.002471576866040492
$ws.Range("D5").Value = 0.1373543769279147
$ws.Range("E5").Value = 0.003262362637362681
$ws.Range("D6").Value = 0.03119840176494664
$ws.Range("E6").Value = -0.004487658937920691
$ws.Range("D7").Value = 0.1198869705409641
$ws.Range("E7").Value = 0.005948446794447948
$ws.Range("D8").Value = 0.1016748525051661
$ws.Range("E8").Value = -0.007756447547023426
$ws.Range("D9").Value = 0.02800155730586086
$ws.Range("E9").Value = -0.008863636363636296
$ws.Range("D10").Value = 0.1246094256935501
$ws.Range("E10").Value = -0.01522130983376724
$ws.Range("D11").Value = 0.2463837561019436
$ws.Range("E11").Value = -0.002405696689761316
$ws.Range("D12").Value = 0.1052169746338834
$ws.Range("E12").Value = 0.006575014943215773
$ws.Range("D13").Value = 0.9999999999999999
$ws.Range("E13").Value = -0.002656903556846779

# Restore sheet protection with the original password.
$ws.Protect("D382")

